# Updating data and plots
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the new Bermuda row (row 30) to the data table.
# Write the alpha_3 code first so the new shared-string entries land in the
# same order as the source file (BMU before Bermuda).
$ws.Range("C30").Value = "BMU"
$ws.Range("A30").Value = "Bermuda"
$ws.Range("B30").Value = "Bermuda"
$ws.Range("D30").Value = 3166
$ws.Range("E30").Value = 142924

# Match the saved view/selection state.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D28").Select()
